$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Update the letter date from September 19 to September 21, 2025.
#    Edit the existing run's paragraph range in place (rather than a
#    Find/Replace) so its original text formatting is preserved as-is.
# ---------------------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -eq "September 19, 2025`r") {
        $para.Range.Text = "September 21, 2025"
        break
    }
}

# ---------------------------------------------------------------------
# 2. Split the mailing-address paragraph
#        "175 Lewis Road Suite, San Jose CA 95111"
#    (the one in the sender/recipient block near the top of the letter,
#    just below "Machining N.M." -- NOT the identical text that also
#    appears later inside the "PROPERTY ADDRESS" table cell, which must
#    stay untouched) into two paragraphs:
#        "175 Lewis Road Suite"
#        "San Jose, CA 95111"
# ---------------------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Information(12)) {
        # Skip any paragraph that lives inside a table (wdWithInTable).
        continue
    }
    if ($para.Range.Text -eq "175 Lewis Road Suite, San Jose CA 95111`r") {
        # Shrink the paragraph down to just the street address, editing
        # the existing run in place so its formatting is kept untouched.
        $para.Range.Text = "175 Lewis Road Suite"

        # Insert a brand-new paragraph right after it (inherits this
        # paragraph's Arial/11pt formatting) and fill it with the
        # city/state/zip.
        $para.Range.InsertParagraphAfter() | Out-Null
        $cityPara = $para.Next()

        # Type the text with a trailing space first, then trim it back
        # off in a second, in-place edit -- this mirrors how Word keeps
        # significant-whitespace markup on a run once it has been
        # through an edit, producing the same "xml:space=preserve" run
        # markup Word itself would leave behind.
        $cityPara.Range.Text = "San Jose, CA 95111 "
        $cityPara.Range.Text = "San Jose, CA 95111"
        break
    }
}

# ---------------------------------------------------------------------
# 3. Remove the empty "NoSpacing" paragraph that immediately follows the
#    "... Board of Directors" line near the signature block.
# ---------------------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t -eq "LEWIS ROAD INDUSTRIAL CENTER Board of Directors`r") {
        $emptyPara = $para.Next()
        if ($emptyPara -ne $null -and $emptyPara.Range.Text -eq "`r") {
            $emptyPara.Range.Delete() | Out-Null
        }
        break
    }
}
